$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "+404.56%"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 22.95
$ws.Range("M2").Value = 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+367.84%"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 22.318
$ws.Range("M3").Value = 3
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "+151.73%"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 25.511
$ws.Range("M4").Value = 3
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "+40.91%"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 12.248
$ws.Range("M5").Value = 3
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "+4.18%"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 6.135
$ws.Range("M6").Value = 3
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "+24.23%"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").Value = 12.391
$ws.Range("M7").Value = 3
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "+67.56%"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").Value = 14.475
$ws.Range("M8").Value = 2
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "+19.93%"
$ws.Range("G9").ClearFormats()
$ws.Range("H9").Value = 22.701
$ws.Range("M9").Value = 3
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "+7.35%"
$ws.Range("G10").ClearFormats()
$ws.Range("H10").Value = 13.169
$ws.Range("M10").Value = 3
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "+2.71%"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = 2.863
$ws.Range("M11").Value = 3
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "+15.62%"
$ws.Range("G12").ClearFormats()
$ws.Range("H12").Value = 9.656000000000001
$ws.Range("M12").Value = 3
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "+182.80%"
$ws.Range("G13").ClearFormats()
$ws.Range("H13").Value = 18.009
$ws.Range("M13").Value = 3
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "+266.74%"
$ws.Range("G14").ClearFormats()
$ws.Range("H14").Value = 17.121
$ws.Range("M14").Value = 3
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "+5.99%"
$ws.Range("G15").ClearFormats()
$ws.Range("H15").Value = 11.68
$ws.Range("M15").Value = 3
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "+51.12%"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").Value = 12.187
$ws.Range("M16").Value = 3

$ws = $wb.Worksheets.Item("Pattern1-Pure Data")
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "+404.56%"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 22.95
$ws.Range("M2").Value = 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+367.84%"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 22.318
$ws.Range("M3").Value = 3
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "+151.73%"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 25.511
$ws.Range("M4").Value = 3
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "+40.91%"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 12.248
$ws.Range("M5").Value = 3
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "+4.18%"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 6.135
$ws.Range("M6").Value = 3

$ws = $wb.Worksheets.Item("Pattern2-Data+Technical")
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "+24.23%"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 12.391
$ws.Range("M2").Value = 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+67.56%"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 14.475
$ws.Range("M3").Value = 2
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "+19.93%"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 22.701
$ws.Range("M4").Value = 3
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "+7.35%"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 13.169
$ws.Range("M5").Value = 3
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "+2.71%"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 2.863
$ws.Range("M6").Value = 3

$ws = $wb.Worksheets.Item("Pattern3-Data+News")
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "+15.62%"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 9.656000000000001
$ws.Range("M2").Value = 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+182.80%"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 18.009
$ws.Range("M3").Value = 3
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "+266.74%"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 17.121
$ws.Range("M4").Value = 3
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "+5.99%"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 11.68
$ws.Range("M5").Value = 3
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "+51.12%"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 12.187
$ws.Range("M6").Value = 3
